$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at the front (shifts B..I <- A..H)
$ws.Columns.Item(1).Insert()

# 2. New "Id" column header + values
$ws.Range("A1").Value = "Id"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A2:A3").NumberFormat = "@"

# 3. Update the (now-shifted) UserName numeric values (D2/D3).
#    D2/D3 are formatted as Text ("@"), so a direct numeric assignment would be
#    coerced to a text value; temporarily switch to a numeric format, set the
#    value, then restore the Text format so the stored cell keeps its original style.
$ws.Range("D2").NumberFormat = "0"
$ws.Range("D2").Value = 1216802002
$ws.Range("D2").NumberFormat = "@"

$ws.Range("D3").NumberFormat = "0"
$ws.Range("D3").Value = 3475716036
$ws.Range("D3").NumberFormat = "@"

# 4. Column width fix-ups (column insert preserves the other widths correctly,
#    only the UserName and Customer columns need new explicit widths).
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666   # -> stored width 11.0
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666   # -> stored width 12.0

# 5. Hyperlinks: Columns.Insert() does not relocate existing Hyperlink anchors,
#    so remove the stale ones and re-create them on the correct (shifted) cells,
#    preserving the original display text/value and the Hyperlink cell style.
$emailStyle = $ws.Range("H2").Style
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:admin@mail.com", [Type]::Missing, [Type]::Missing, "mailto:admin@mail.com")
$ws.Range("H2").Value = "admin@mail.com"
$ws.Range("H2").Style = $emailStyle

$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:customer@mail.com", [Type]::Missing, [Type]::Missing, "mailto:customer@mail.com")
$ws.Range("H3").Value = "customer@mail.com"
$ws.Range("H3").Style = $emailStyle

# 6. Update the active selection to match the authored workbook.
$ws.Range("E9").Select()

Write-Host "done"
